# Insert a new data row at row 23 (pushes old rows 23..120 down to 24..121,
# matching the rest of the table's existing formatting/styles automatically),
# then populate the newly inserted row with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(23).Insert()

$ws.Range("A23").Value = 3
$ws.Range("B23").Value = "Femacal de La Calera"
$ws.Range("C23").Value = "Coquimbo"
$ws.Range("D23").Value = 44558
$ws.Range("E23").Value = 5
$ws.Range("F23").Value = 100112030
$ws.Range("G23").Value = "Poroto granado"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 73
$ws.Range("K23").Value = 41000
$ws.Range("L23").Value = 42000
$ws.Range("M23").Value = 41521
$ws.Range("N23").Value = "`$/saco 25 kilos"
$ws.Range("O23").Value = "Provincia de Talca"
$ws.Range("P23").Value = 1661
$ws.Range("Q23").Value = 25
$ws.Range("R23").Value = "Hortaliza"
